$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("constants")

# Remove the three organ-specific case-proportion rows (rows 3-5), which were
# previously computed via formulas from raw notification counts. They are
# replaced elsewhere in the sheet with fixed proportions taken from the GTB
# report, so these rows are deleted outright (everything below shifts up).
$ws.Range("A3:E5").EntireRow.Delete() | Out-Null
$ws.Range("B4").Select() | Out-Null
